$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.272.63'
$ws.Range("E2").Value = '  +1.05%  '

$ws.Range("D3").Value = '1.884.35'
$ws.Range("E3").Value = '  +1.38%  '

$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = '  +0.33%  '

$ws.Range("D5").Value = "'314.44"
$ws.Range("E5").Value = '  +0.92%  '

$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = '  +0.39%  '

$ws.Range("D7").Value = "'0.5138"
$ws.Range("E7").Value = '  +0.18%  '

$ws.Range("D8").Value = "'0.3919"
$ws.Range("E8").Value = '  +3.13%  '

$ws.Range("D9").Value = "'0.08367"
$ws.Range("E9").Value = '  +0.71%  '

$ws.Range("D10").Value = "'1.124"
$ws.Range("E10").Value = '  +1.64%  '

$ws.Range("D11").Value = "'41.62"
$ws.Range("E11").Value = '  +0.82%  '

$ws.Range("D12").Value = "'6.236"
$ws.Range("E12").Value = '  +0.76%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.894.63'
$ws.Range("E13").Value = '  +1.75%  '

$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").Value = "'20.74"
$ws.Range("E14").Value = '  +1.62%  '

$ws.Range("D15").Value = "'7.293"
$ws.Range("E15").Value = '  +1.70%  '

$ws.Range("E16").Value = '  +0.30%  '

$ws.Range("D17").Value = "'0.00001107"
$ws.Range("E17").Value = '  +1.35%  '

$ws.Range("D18").Value = "'91.54"
$ws.Range("E18").Value = '  +1.35%  '

$ws.Range("D19").Value = "'0.06670"
$ws.Range("E19").Value = '  +0.88%  '

$ws.Range("D20").Value = "'17.84"
$ws.Range("E20").Value = '  +0.40%  '

$ws.Range("D21").Value = "'1.006"
$ws.Range("E21").Value = '  +0.40%  '

$ws.Range("D22").Value = "'6.053"

$ws.Range("D23").Value = '28.313.43'
$ws.Range("E23").Value = '  +1.08%  '

$ws.Range("D24").Value = "'11.18"
$ws.Range("E24").Value = '  +1.23%  '

$ws.Range("D25").Value = "'2.280"
$ws.Range("E25").Value = '  +1.33%  '

$ws.Range("D26").Value = '2.095.78'
$ws.Range("E26").Value = '  +0.90%  '

$ws.Range("D27").Value = "'2.526"
$ws.Range("E27").Value = '  -1.58%  '

$ws.Range("D28").Value = "'159.23"
$ws.Range("E28").Value = '  +1.28%  '

$ws.Range("D29").Value = "'20.68"
$ws.Range("E29").Value = '  +1.45%  '

$ws.Range("D30").Value = "'125.75"
$ws.Range("E30").Value = '  +0.39%  '

$ws.Range("D31").Value = "'0.1068"
$ws.Range("E31").Value = '  +0.88%  '

$ws.Range("D32").Value = "'1.052"
$ws.Range("E32").Value = '  +1.28%  '

$ws.Range("D33").Value = "'5.900"
$ws.Range("E33").Value = '  +5.66%  '

$ws.Range("D34").Value = "'3.607"
$ws.Range("E34").Value = '  +0.12%  '

$ws.Range("D35").Value = "'9.764"
$ws.Range("E35").Value = '  +2.07%  '

$ws.Range("E36").Value = '  +2.53%  '

$ws.Range("D37").Value = "'0.06597"
$ws.Range("E37").Value = '  +1.31%  '

$ws.Range("D38").Value = "'0.2198"
$ws.Range("E38").Value = '  +2.26%  '

$ws.Range("D39").Value = "'1.214"
$ws.Range("E39").Value = '  +0.63%  '

$ws.Range("D40").Value = "'0.6542"
$ws.Range("E40").Value = '  +2.33%  '

$ws.Range("D41").Value = "'5.027"
$ws.Range("E41").Value = '  +3.48%  '

$ws.Range("D42").Value = "'1.231"
$ws.Range("E42").Value = '  +0.38%  '

$ws.Range("D43").Value = "'11.32"
$ws.Range("E43").Value = '  +0.65%  '

$ws.Range("D44").Value = "'0.6170"
$ws.Range("E44").Value = '  +1.86%  '

$ws.Range("D45").Value = "'13.15"
$ws.Range("E45").Value = '  +0.54%  '

$ws.Range("D46").Value = "'1.288"
$ws.Range("E46").Value = '  +0.37%  '

$ws.Range("D47").Value = "'3.687"
$ws.Range("E47").Value = '  +0.88%  '

$ws.Range("D48").Value = "'2.021"
$ws.Range("E48").Value = '  +2.46%  '

$ws.Range("D49").Value = "'1.239"
$ws.Range("E49").Value = '  +2.70%  '

$ws.Range("D50").Value = "'121.73"
$ws.Range("E50").Value = '  +1.17%  '

$ws.Range("D51").Value = "'79.25"
$ws.Range("E51").Value = '  -0.35%  '
